$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "AU-4,AU-14 (1)"
$ws.Range("A4").Value = "CM-6 b,AU-4"
$ws.Range("A5").Value = "CM-6 b,SC-5 (2),SC-5"
$ws.Range("A6").Value = "CM-5 (1),AC-6 (8),AU-8 b,AU-7 b,AU-7 a,AU-12 (3),AC-6 (9)"
$ws.Range("C6").Value = "SRG-OS-000326-GPOS-00126,SRG-OS-000327-GPOS-00127,SRG-OS-000337-GPOS-00129,SRG-OS-000348-GPOS-00136,SRG-OS-000349-GPOS-00137,SRG-OS-000350-GPOS-00138,SRG-OS-000351-GPOS-00139,SRG-OS-000352-GPOS-00140,SRG-OS-000353-GPOS-00141,SRG-OS-000354-GPOS-00142,SRG-OS-000358-GPOS-00145,SRG-OS-000365-GPOS-00152"
$ws.Range("A7").Value = "CM-6 b,AU-12 c,CM-5 (1),AU-8 b,AU-7 b,AU-7 a,AU-12 (3),AU-12 a"
$ws.Range("K7").Value = "Run the following command to determine if the  audit  package is installed:  `$ rpm -q audit `n`nIf the audit package is not installed then this is a finding."
$ws.Range("A14").Value = "CM-7 (5) (b),CM-7 (2)"
$ws.Range("A15").Value = "CM-7 (5) (b),CM-7 (2)"
$ws.Range("A17").Value = "CM-6 b,CM-7 (2)"
$ws.Range("A22").Value = "CM-6 b,CM-7 (2)"
$ws.Range("A23").Value = "CM-6 b,CM-7 (2)"
$ws.Range("A38").Value = "AC-7 b,AC-7 a"
$ws.Range("A39").Value = "AC-7 b,AC-7 a"
$ws.Range("A40").Value = "AC-7 b,AC-7 a"
$ws.Range("A41").Value = "AC-7 b,AC-7 a"
$ws.Range("A45").Value = "AU-3 (1),IA-2,IA-8"
$ws.Range("A46").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A47").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A48").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A49").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A50").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A51").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A52").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A53").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A54").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A55").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A56").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A57").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A58").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A59").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A60").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A61").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A62").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A63").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A64").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A65").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A66").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A67").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A68").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A69").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A70").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A71").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A72").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A73").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A74").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A75").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A76").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A77").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A78").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A79").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("M79").Value = "Configure Red Hat Enterprise Linux 9 to generate audit records for all account creations, modifications, disabling, and termination events that affect  /var/log/lastlog .`nAdd or update the following file system rule to  /etc/audit/rules.d/audit.rules :`n-w /var/log/lastlog -p wa -k logins`nThe audit daemon must be restarted for the changes to take effect."
$ws.Range("A80").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A81").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A82").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A83").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A84").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A85").Value = "AU-3 (1),AU-3,MA-4 (1) (a)"
$ws.Range("A86").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A87").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A88").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A89").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A90").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A91").Value = "AU-3 (1),AU-12 c,MA-4 (1) (a)"
$ws.Range("A92").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A93").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A94").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A95").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A96").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A97").Value = "AU-3 (1),AU-12 c,AU-3,MA-4 (1) (a)"
$ws.Range("A98").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A99").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A100").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A101").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A102").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-12 a"
$ws.Range("A103").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AC-2 (4),AU-12 a"
$ws.Range("A104").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AC-2 (4),AU-12 a"
$ws.Range("A105").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AC-2 (4),AU-12 a"
$ws.Range("A106").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AC-2 (4)"
$ws.Range("A107").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AC-2 (4),AU-12 a"
$ws.Range("A108").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AC-2 (4),AU-12 a"
$ws.Range("A109").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AC-2 (4),AU-12 a"
$ws.Range("A110").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AC-2 (4),AU-12 a"
$ws.Range("A111").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AC-2 (4),AU-12 a"
$ws.Range("A112").Value = "AU-3,AU-12 c,AU-3 (1),MA-4 (1) (a),AU-14 (1),AU-12 a"
$ws.Range("M112").Value = "To ensure that  audit=1  is added as a kernel command line`nargument to newly installed kernels, add  audit=1  to the`ndefault Grub2 command line for Linux operating systems.  Modify the line within`n /etc/default/grub  as shown below:`n GRUB_CMDLINE_LINUX=`"... audit=1 ...`" `nRun the following command to update command line for already installed kernels: # grubby --update-kernel=ALL --args=`"audit=1`" "
$ws.Range("A113").Value = "CM-6 b,AC-6 (10)"
$ws.Range("A114").Value = "CM-6 b,AC-6 (10)"
$ws.Range("A115").Value = "AC-11 b,AC-6 (10)"
$ws.Range("A116").Value = "CM-6 b,AC-6 (10)"
$ws.Range("A120").Value = "AU-12 c,AU-3,MA-4 (1) (a),AU-12 a"
$ws.Range("A121").Value = "AU-12 c,AU-3,MA-4 (1) (a),AU-12 a"
$ws.Range("A127").Value = "CM-5 (1),AU-12 c,AC-2 (4),AC-6 (9)"
$ws.Range("A133").Value = "SC-8,AC-17 (2),MA-4 c,SC-13"
$ws.Range("A134").Value = "MA-4 (7),MA-4 e,AC-12,SC-10"
$ws.Range("A138").Value = "AU-3,AU-6 (4),AU-7 (1),CM-6 b,AU-14 (1),AU-3 (1),CM-5 (1),MA-4 (1) (a),AU-7 a,AU-12 a"
$ws.Range("A141").Value = "AU-9 (3),AU-9"
$ws.Range("A142").Value = "AU-9 (3),AU-9"
$ws.Range("A143").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A144").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A145").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A146").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A147").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A148").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A149").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A150").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A151").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A152").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A153").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A154").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A155").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A166").Value = "SC-8,SC-8 (1),SC-8 (2)"
$ws.Range("A179").Value = "CM-6 b,AU-4 (1),AU-6 (4)"
$ws.Range("A180").Value = "AC-17 (1),CM-7 b,CM-6 b,AC-17 (9)"
$ws.Range("A181").Value = "AC-17 (1),CM-7 b,CM-6 b"
$ws.Range("A182").Value = "AU-9,SI-11 b"
$ws.Range("A183").Value = "AU-9,SI-11 b"
$ws.Range("A184").Value = "AU-9,SI-11 b"
$ws.Range("A185").Value = "AU-9,SI-11 b"
$ws.Range("A186").Value = "AU-9,SI-11 b"
$ws.Range("A187").Value = "AU-9,SI-11 b"
$ws.Range("A194").Value = "CM-6 b,AU-3"
$ws.Range("A207").Value = "CM-6 b,AU-4 (1)"
$ws.Range("A216").Value = "AU-12 c,AC-2 (4),AC-6 (9)"
$ws.Range("A221").Value = "CM-6 b,IA-2 (5)"
$ws.Range("A222").Value = "IA-2 (5),IA-2,IA-2 (4),IA-2 (2),IA-2 (3)"
$ws.Range("A223").Value = "IA-2 (5),IA-2,IA-2 (4),IA-2 (2),IA-2 (3)"
$ws.Range("A224").Value = "AC-18 (1),SC-8 (1),SC-8"
$ws.Range("A227").Value = "CM-6 b,IA-7"
$ws.Range("A228").Value = "CM-6 b,IA-7"
$ws.Range("A229").Value = "CM-6 b,IA-7"
$ws.Range("A232").Value = "MA-4 (6),SC-13"
$ws.Range("A233").Value = "AC-17 (2),MA-4 (6)"
$ws.Range("A234").Value = "MA-4 (6),SC-13"
$ws.Range("A242").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("K242").Value = "To verify that auditing is configured for system administrator actions, run the following command:`n `$ sudo auditctl -l | grep `"watch=/var/run/faillock\|-w /var/run/faillock`" `n`nIf there is no output then this is a finding."
$ws.Range("M242").Value = "Configure Red Hat Enterprise Linux 9 to generate audit records for all account creations, modifications, disabling, and termination events that affect  /var/run/faillock .`nAdd or update the following file system rule to  /etc/audit/rules.d/audit.rules :`n-w /var/run/faillock -p wa -k logins`nThe audit daemon must be restarted for the changes to take effect."
$ws.Range("A243").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("M243").Value = "Configure Red Hat Enterprise Linux 9 to generate audit records for all account creations, modifications, disabling, and termination events that affect  /var/log/tallylog .`nAdd or update the following file system rule to  /etc/audit/rules.d/audit.rules :`n-w /var/log/tallylog -p wa -k logins`nThe audit daemon must be restarted for the changes to take effect."
$ws.Range("A245").Value = "CM-6 b,SI-16,SC-2"
$ws.Range("A271").Value = "IA-2 (2),IA-2 (3),IA-2 (1),IA-2 (4)"
$ws.Range("A273").Value = "AU-12 c,MA-4 (1) (a)"
$ws.Range("A276").Value = "CM-6 b,SC-4"
$ws.Range("A277").Value = "SC-4,SC-2"
$ws.Range("A278").Value = "SC-4,SC-2"
$ws.Range("A281").Value = "CM-6 b,AU-12 a"
$ws.Range("A284").Value = "CM-6 b,CM-5 (3)"
$ws.Range("A310").Value = "AU-8 (1) (b),AU-8 (1) (a),AU-8 b"
$ws.Range("A342").Value = "IA-3,CM-7 b"
$ws.Range("A345").Value = "AC-17 (1),CM-7 b"
$ws.Range("A347").Value = "CM-6 b,IA-5 (1) (c),CM-7 a"
$ws.Range("A358").Value = "AC-11 (1),AC-11 b"
$ws.Range("A361").Value = "CM-3 (5),SI-6 b,SI-6 d"
$ws.Range("K361").Value = "Verify the operating system routinely checks the baseline configuration for unauthorized changes.`n`nTo determine that periodic AIDE execution has been scheduled, run the following command:`n `$ grep aide /etc/crontab `nThe output should return something similar to the following:`n`n 05 4 * * * root  --check `n`n`nNOTE: The usage of special cron times, such as @daily or @weekly, is acceptable.`n`nIf AIDE is not configured to scan periodically then this is a finding."
$ws.Range("A362").Value = "CM-6 b,CM-7 a"
$ws.Range("M367").Value = "To ensure that  pti=on  is added as a kernel command line`nargument to newly installed kernels, add  pti=on  to the`ndefault Grub2 command line for Linux operating systems.  Modify the line within`n /etc/default/grub  as shown below:`n GRUB_CMDLINE_LINUX=`"... pti=on ...`" `nRun the following command to update command line for already installed kernels: # grubby --update-kernel=ALL --args=`"pti=on`" "
$ws.Range("A374").Value = "CM-6 b,CM-7 a"
$ws.Range("A375").Value = "CM-6 b,CM-7 a"
$ws.Range("A376").Value = "CM-6 b,CM-7 a"
$ws.Range("A385").Value = "CM-6 b,AC-17 (2)"
$ws.Range("A391").Value = "CM-6 b,IA-5 (1) (a)"
$ws.Range("A398").Value = "CM-6 b,SI-16"
$ws.Range("M398").Value = "Add or edit the following line in a system configuration file in the `"/etc/sysctl.d/`" directory:`nkernel.randomize_va_space = 2`nLoad settings from all system configuration files with the following command:`n`$ sudo sysctl --system"
$ws.Range("F399").Value = "Red Hat Enterprise Linux 9 must, for networked systems, compare internal information system clocks at least every 24 hours with a server which is synchronized to one of the redundant United States Naval Observatory (USNO) time servers, or a time server designated for the appropriate DoD network (NIPRNet/SIPRNet), and/or the Global Positioning System (GPS)."
$ws.Range("A401").Value = "CM-6 b,SC-3"
$ws.Range("A402").Value = "CM-6 b,SC-3"
$ws.Range("A403").Value = "CM-6 b,SC-3"
$ws.Range("A448").Value = "CM-6 b,IA-5 (1) (c)"
$ws.Range("A450").Value = "CM-6 b,CM-5 (1)"
$ws.Range("A451").Value = "CM-6 b,CM-5 (1)"
$ws.Range("A524").Value = "CM-6 b,SC-2"
$ws.Range("A525").Value = "CM-6 b,SC-2"
$ws.Range("A541").Value = "CM-6 b,SI-2 (2)"
$ws.Range("A550").Value = "CM-6 b,SI-2 (2)"
$ws.Range("A558").Value = "CM-3 (5),SI-6 a"
